# Fixes #67 - Implemented blank node with edit form support for related resource.
# Adds a new "related resource" column (T) to the import-test fixture sheet,
# mirroring the bold header style used by the other header cells and the
# selection/scroll state left behind after entering the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell: "related resource" - use the same bold header style (s="2")
# as the other column headers in row 1 (e.g. C1, D1, ...).
$ws.Range("T1").Value = "related resource"
$ws.Range("T1").Font.Bold = $true

# New data cell for the first data row.
$ws.Range("T2").Value = "RELATED RESOURCE @{relatedType=relation; url=http://test.com/related_resource/relation}"

# Leave the selection on the newly entered cell, matching the workbook state
# after the edit was made.
$ws.Range("T2").Select()
